$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the picture shape that needs to move (cNvPr id="852", name "Picture 2").
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 852) {
        $target = $sh
        break
    }
}

if ($target -ne $null) {
    # Move it from near the end of the z-order to just after the first
    # shape (the other "Picture 2", id 1026) near the front of the z-order:
    # send it all the way to the back, then bring it forward one step.
    $target.ZOrder(1)   # msoSendToBack
    $target.ZOrder(2)   # msoBringForward
}
